$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as text in this sheet
# (values like "600.00" / "69.008.64" are not valid numbers). Force text
# format before writing so Excel does not auto-coerce/trim them as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.006.75"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.796.59"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "600.15"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "163.79"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").Value = "3.795.66"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "36.98"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "4.427.75"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "3.811.16"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "69.137.66"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "17.21"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "11.38"
$ws.Range("E21").Value = "  +6.12%  "
$ws.Range("D22").Value = "485.51"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "0.0000159"
$ws.Range("E24").Value = "  +6.18%  "
$ws.Range("D25").Value = "84.42"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -4.60%  "
$ws.Range("D33").Value = "3.954.99"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").Value = "31.55"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "3.740.28"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("D39").Value = "5.83"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("D42").Value = "0.318"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "436.30"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "48.46"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D47").Value = "8.35"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "2.820.43"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "141.55"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.29"
$ws.Range("E50").Value = "  +5.56%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "39.13"
$ws.Range("E51").Value = "  -2.46%  "
